# Scheduled market-price refresh: update currentAveragePrice / NQ / HQ price
# columns (H:N) for the rows whose crafted-item market data changed, across
# the ALC, BSM, CRP, CUL, GSM, LTW and WVR sheets. Values only - no
# structural changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2171.7144
$ws.Range("I40").Value = 1633.3334
$ws.Range("K40").Value = 1633.3334
$ws.Range("M40").Value = -1458.3334
$ws.Range("H43").Value = 1994.2142
$ws.Range("I43").Value = 4460.4
$ws.Range("J43").Value = 624.1111
$ws.Range("K43").Value = 4460.4
$ws.Range("L43").Value = 624.1111
$ws.Range("M43").Value = -4391.4
$ws.Range("N43").Value = -762.1111
$ws.Range("H88").Value = 1060334.8
$ws.Range("I88").Value = 2629
$ws.Range("J88").Value = 1985827.4
$ws.Range("K88").Value = 2629
$ws.Range("L88").Value = 1985827.4
$ws.Range("M88").Value = -2223
$ws.Range("N88").Value = -1986639.4
$ws.Range("H91").Value = 1060334.8
$ws.Range("I91").Value = 2629
$ws.Range("J91").Value = 1985827.4
$ws.Range("K91").Value = 2629
$ws.Range("L91").Value = 1985827.4
$ws.Range("M91").Value = -1225
$ws.Range("N91").Value = -1988635.4
$ws.Range("H98").Value = 1423.2307
$ws.Range("I98").Value = 1530.3478
$ws.Range("J98").Value = 602
$ws.Range("K98").Value = 1530.3478
$ws.Range("L98").Value = 602
$ws.Range("M98").Value = -32.34780000000001
$ws.Range("N98").Value = -3598
$ws.Range("H100").Value = 2059.0908
$ws.Range("I100").Value = 3000
$ws.Range("J100").Value = 1850
$ws.Range("K100").Value = 3000
$ws.Range("L100").Value = 1850
$ws.Range("M100").Value = -2459
$ws.Range("N100").Value = -2932
$ws.Range("H103").Value = 953.6
$ws.Range("I103").Value = 1370.5714
$ws.Range("J103").Value = 588.75
$ws.Range("K103").Value = 4111.7142
$ws.Range("L103").Value = 1766.25
$ws.Range("M103").Value = -3525.7142
$ws.Range("N103").Value = -2938.25
$ws.Range("H116").Value = 3679.7334
$ws.Range("I116").Value = 3608.182
$ws.Range("J116").Value = 3876.5
$ws.Range("K116").Value = 3608.182
$ws.Range("L116").Value = 3876.5
$ws.Range("M116").Value = -166.1819999999998
$ws.Range("N116").Value = -10760.5
$ws.Range("H122").Value = 1423.2307
$ws.Range("I122").Value = 1530.3478
$ws.Range("J122").Value = 602
$ws.Range("K122").Value = 4591.0434
$ws.Range("L122").Value = 1806
$ws.Range("M122").Value = -2141.0434
$ws.Range("N122").Value = -6706
$ws.Range("H125").Value = 7124.9
$ws.Range("I125").Value = 2255.1667
$ws.Range("J125").Value = 14429.5
$ws.Range("K125").Value = 20296.5003
$ws.Range("L125").Value = 129865.5
$ws.Range("M125").Value = -17836.5003
$ws.Range("N125").Value = -134785.5
$ws.Range("H132").Value = 2036.9788
$ws.Range("I132").Value = 1933.45
$ws.Range("J132").Value = 2628.5715
$ws.Range("K132").Value = 5800.35
$ws.Range("L132").Value = 7885.7145
$ws.Range("M132").Value = -3270.35
$ws.Range("N132").Value = -12945.7145
$ws.Range("H135").Value = 1736.5714
$ws.Range("I135").Value = 1023.2941
$ws.Range("J135").Value = 4768
$ws.Range("K135").Value = 9209.6469
$ws.Range("L135").Value = 42912
$ws.Range("M135").Value = -6674.6469
$ws.Range("N135").Value = -47982
$ws.Range("H137").Value = 939.5833
$ws.Range("I137").Value = 827.55
$ws.Range("K137").Value = 2482.65
$ws.Range("M137").Value = 67.35000000000036
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 50034084
$ws.Range("I20").Value = 41360.5
$ws.Range("J20").Value = 250004980
$ws.Range("K20").Value = 41360.5
$ws.Range("L20").Value = 250004980
$ws.Range("M20").Value = -41113.5
$ws.Range("N20").Value = -250005474
$ws.Range("H64").Value = 1585.9231
$ws.Range("I64").Value = 2681.2
$ws.Range("J64").Value = 901.375
$ws.Range("K64").Value = 2681.2
$ws.Range("L64").Value = 901.375
$ws.Range("M64").Value = -2456.2
$ws.Range("N64").Value = -1351.375
$ws.Range("H67").Value = 1585.9231
$ws.Range("I67").Value = 2681.2
$ws.Range("J67").Value = 901.375
$ws.Range("K67").Value = 2681.2
$ws.Range("L67").Value = 901.375
$ws.Range("M67").Value = -1901.2
$ws.Range("N67").Value = -2461.375
$ws.Range("H80").Value = 566.82355
$ws.Range("I80").Value = 283.7143
$ws.Range("J80").Value = 765
$ws.Range("K80").Value = 283.7143
$ws.Range("L80").Value = 765
$ws.Range("M80").Value = 714.2857
$ws.Range("N80").Value = -2761
$ws.Range("H83").Value = 566.82355
$ws.Range("I83").Value = 283.7143
$ws.Range("J83").Value = 765
$ws.Range("K83").Value = 1418.5715
$ws.Range("L83").Value = 3825
$ws.Range("M83").Value = 3573.4285
$ws.Range("N83").Value = -13809
$ws.Range("H94").Value = 2463.8
$ws.Range("I94").Value = 1769.6666
$ws.Range("J94").Value = 3505
$ws.Range("K94").Value = 1769.6666
$ws.Range("L94").Value = 3505
$ws.Range("M94").Value = -1318.6666
$ws.Range("N94").Value = -4407
$ws.Range("H99").Value = 47621224
$ws.Range("I99").Value = 71430900
$ws.Range("J99").Value = 1885.7142
$ws.Range("K99").Value = 71430900
$ws.Range("L99").Value = 1885.7142
$ws.Range("M99").Value = -71429402
$ws.Range("N99").Value = -4881.7142
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2076.1155
$ws.Range("I99").Value = 1522.8823
$ws.Range("J99").Value = 3121.111
$ws.Range("K99").Value = 1522.8823
$ws.Range("L99").Value = 3121.111
$ws.Range("M99").Value = -24.88229999999999
$ws.Range("N99").Value = -6117.111
$ws.Range("H126").Value = 2076.1155
$ws.Range("I126").Value = 1522.8823
$ws.Range("J126").Value = 3121.111
$ws.Range("K126").Value = 4568.6469
$ws.Range("L126").Value = 9363.332999999999
$ws.Range("M126").Value = -2098.6469
$ws.Range("N126").Value = -14303.333
$ws.Range("H134").Value = 1966.2812
$ws.Range("I134").Value = 1235.4231
$ws.Range("J134").Value = 5133.3335
$ws.Range("K134").Value = 3706.2693
$ws.Range("L134").Value = 15400.0005
$ws.Range("M134").Value = -1171.2693
$ws.Range("N134").Value = -20470.0005
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 345847.3
$ws.Range("I122").Value = 814.36365
$ws.Range("J122").Value = 556700.75
$ws.Range("K122").Value = 7329.27285
$ws.Range("L122").Value = 5010306.75
$ws.Range("M122").Value = -4879.27285
$ws.Range("N122").Value = -5015206.75
$ws.Range("H132").Value = 660
$ws.Range("I132").Value = 546.6667
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 4920.0003
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -2390.0003
$ws.Range("N132").Value = -14060
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 734716.1
$ws.Range("I122").Value = 2195482.2
$ws.Range("J122").Value = 4333
$ws.Range("K122").Value = 6586446.600000001
$ws.Range("L122").Value = 12999
$ws.Range("M122").Value = -6583996.600000001
$ws.Range("N122").Value = -17899
$ws.Range("H126").Value = 3678.111
$ws.Range("I126").Value = 3781.8462
$ws.Range("J126").Value = 3408.4
$ws.Range("K126").Value = 11345.5386
$ws.Range("L126").Value = 10225.2
$ws.Range("M126").Value = -8875.5386
$ws.Range("N126").Value = -15165.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2051.5
$ws.Range("I7").Value = 2106.1667
$ws.Range("J7").Value = 1887.5
$ws.Range("K7").Value = 2106.1667
$ws.Range("L7").Value = 1887.5
$ws.Range("M7").Value = -1994.1667
$ws.Range("N7").Value = -2111.5
$ws.Range("H16").Value = 612.13336
$ws.Range("I16").Value = 532.44446
$ws.Range("J16").Value = 731.6667
$ws.Range("K16").Value = 532.44446
$ws.Range("L16").Value = 731.6667
$ws.Range("M16").Value = -362.44446
$ws.Range("N16").Value = -1071.6667
$ws.Range("H46").Value = 1366.6666
$ws.Range("I46").Value = 1633.3334
$ws.Range("J46").Value = 833.3333
$ws.Range("K46").Value = 1633.3334
$ws.Range("L46").Value = 833.3333
$ws.Range("M46").Value = -1445.3334
$ws.Range("N46").Value = -1209.3333
$ws.Range("H55").Value = 286.14285
$ws.Range("I55").Value = 220.5
$ws.Range("J55").Value = 450.25
$ws.Range("K55").Value = 220.5
$ws.Range("L55").Value = 450.25
$ws.Range("M55").Value = -47.5
$ws.Range("N55").Value = -796.25
$ws.Range("H93").Value = 1578.4166
$ws.Range("I93").Value = 1478.1428
$ws.Range("K93").Value = 1478.1428
$ws.Range("M93").Value = -230.1428000000001
$ws.Range("H122").Value = 3726.3962
$ws.Range("I122").Value = 4884.6294
$ws.Range("J122").Value = 2523.6155
$ws.Range("K122").Value = 14653.8882
$ws.Range("L122").Value = 7570.8465
$ws.Range("M122").Value = -12203.8882
$ws.Range("N122").Value = -12470.8465
$ws.Range("H126").Value = 2051.5
$ws.Range("I126").Value = 2106.1667
$ws.Range("J126").Value = 1887.5
$ws.Range("K126").Value = 6318.500100000001
$ws.Range("L126").Value = 5662.5
$ws.Range("M126").Value = -3848.500100000001
$ws.Range("N126").Value = -10602.5
$ws.Range("H132").Value = 3439.5625
$ws.Range("I132").Value = 3068.12
$ws.Range("J132").Value = 4766.143
$ws.Range("K132").Value = 9204.360000000001
$ws.Range("L132").Value = 14298.429
$ws.Range("M132").Value = -6674.360000000001
$ws.Range("N132").Value = -19358.429
$ws.Range("H136").Value = 5003.609
$ws.Range("I136").Value = 1454.15
$ws.Range("J136").Value = 28666.666
$ws.Range("K136").Value = 4362.450000000001
$ws.Range("L136").Value = 85999.99800000001
$ws.Range("M136").Value = -1812.450000000001
$ws.Range("N136").Value = -91099.99800000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 77728.75
$ws.Range("J133").Value = 77728.75
$ws.Range("L133").Value = 77728.75
$ws.Range("N133").Value = -87848.75
